$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.003994804209775715
$ws.Range("C2").Value = 0.00000000222304730179701
$ws.Range("D2").Value = 26.21740644021617
$ws.Range("E2").Value = 0.496779210170732
$ws.Range("G2").Value = 26.71818045681972

$ws.Range("B3").Value = 0.3048080303191223
$ws.Range("C3").Value = 0.3127903958511391
$ws.Range("D3").Value = 26.21740644021617
$ws.Range("E3").Value = 645.3272768299601
$ws.Range("G3").Value = 672.1622816963464

$ws.Range("B4").Value = 3.230985683306322
$ws.Range("C4").Value = 1.667794583268128
$ws.Range("D4").Value = 0.8054896365839992
$ws.Range("E4").Value = 0.496779210170732
$ws.Range("G4").Value = 6.201049113329182

$ws.Range("B5").Value = 0.127881588408715
$ws.Range("C5").Value = 0.3127903958511391
$ws.Range("D5").Value = 0.1575252929769615
$ws.Range("E5").Value = 0.496779210170732
$ws.Range("G5").Value = 1.094976487407548
